$wb = $excel.ActiveWorkbook

# --- Logs sheet: append new row 59 -----------------------------------
$ws = $wb.Worksheets.Item("Logs")

$ws.Range("A59").Value = "Kun je 200 M8-bouten RVS bestellen voor Van Dijk?"
$ws.Range("B59").Value = "mailmind.test@zohomail.eu"
$ws.Range("C59").Value = "Testmail #18: Kun je 200 M8-bouten RVS bestellen voor Van Dijk?"
$ws.Range("D59").Value = "Inkoop / Bestellingen"
$ws.Range("E59").Value = "Beste klant,`nBedankt voor uw e-mail. Om uw verzoek te kunnen verwerken, hebben we wat meer informatie nodig. Kunt u ons misschien laten weten van welke leverancier u de M8-bouten RVS wilt bestellen en wat het gewenste aantal is? Op die manier kunnen we uw bestelling nauwkeurig verwerken.`nMet vriendelijke groet,`n[Naam bedrijf]"
$ws.Range("F59").Value = "2025-08-05 20:00:40"
$ws.Range("G59").Value = "Ja"
$ws.Range("H59").Value = "Nee"
$ws.Range("I59").Value = "Ja"
$ws.Range("J59").Value = "Nee"

# writing the multi-line E59 text auto-bumps the row height; put it back
# to the sheet's normal (auto) height, same as every other data row.
$ws.Rows(59).EntireRow.AutoFit()

# --- extend conditional formatting ranges to cover the new row -------
$ws.Range("D2:D58").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D59"))
$ws.Range("G2:G58").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G59"))
$ws.Range("H2:H58").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H59"))
$ws.Range("I2:I58").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I59"))
$ws.Range("J2:J58").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J2:J59"))

# --- Dashboard sheet: bump Inkoop / Bestellingen count ----------------
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B4").Value = 7
